# Apply the workbook edit described by the target diff:
#  1. Update the "as of" date in the confidential disclosure note (A41)
#     from 2021-05-25 to 2021-05-26.
#  2. Refresh the "Weight" (column D) and "Percent Change" (column E)
#     values for the holdings rows (rows 2-38) to the newer snapshot.
#
# The worksheet is protected (no-password legacy protection), so it must
# be unprotected before the cell writes land, then protection is restored
# afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect("")

# --- Confidential footer date refresh (A41) ---
$ws.Range("A41").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-26 for illustrative purposes only and are subject to change."

# --- Weight (D) / Percent Change (E) refresh for rows 2-38 ---
$ws.Range("D2").Value = 0.03113897503333115
$ws.Range("E2").Value = 0.001222067039106101
$ws.Range("D3").Value = 0.02941878275233483
$ws.Range("E3").Value = 0.003761780312029783
$ws.Range("D4").Value = 0.02872177481906431
$ws.Range("E4").Value = 0.006164911379398896
$ws.Range("D5").Value = 0.0632753121149104
$ws.Range("E5").Value = 0.001874779460272036
$ws.Range("D6").Value = 0.01538853726769395
$ws.Range("E6").Value = -0.008099924299772687
$ws.Range("D7").Value = 0.01558074836906939
$ws.Range("E7").Value = -0.004934579439252351
$ws.Range("D8").Value = 0.02819096153506387
$ws.Range("E8").Value = -0.01425619834710734
$ws.Range("D9").Value = 0.03533966128520885
$ws.Range("E9").Value = 0.01219646192726076
$ws.Range("D10").Value = 0.02920094350410933
$ws.Range("E10").Value = 0.002034547412933296
$ws.Range("D11").Value = 0.02869401099331008
$ws.Range("E11").Value = 0.003078671908303132
$ws.Range("D12").Value = 0.01104437222398146
$ws.Range("E12").Value = 0.02241364155752823
$ws.Range("D13").Value = 0.01440612497177506
$ws.Range("E13").Value = -0.0007547169811320531
$ws.Range("D14").Value = 0.01414382477181727
$ws.Range("E14").Value = 0.005092726049774177
$ws.Range("D15").Value = 0.008876658136247382
$ws.Range("E15").Value = 0.01014873140857397
$ws.Range("D16").Value = 0.008083350499661494
$ws.Range("E16").Value = 0.0105923043666234
$ws.Range("D17").Value = 0.03094734638983867
$ws.Range("E17").Value = 0.00865135479337753
$ws.Range("D18").Value = 0.02460942801276825
$ws.Range("E18").Value = -0.003471318233099097
$ws.Range("D19").Value = 0.0335379249005987
$ws.Range("E19").Value = -0.00521014241055906
$ws.Range("D20").Value = 0.03182064490901716
$ws.Range("E20").Value = -0.0003965953811890799
$ws.Range("D21").Value = 0.04677272706974953
$ws.Range("E21").Value = 0.01015329567011336
$ws.Range("D22").Value = 0.03538528715270706
$ws.Range("E22").Value = 0.007928451894323896
$ws.Range("D23").Value = 0.03074892240437837
$ws.Range("E23").Value = 0.001831097079715738
$ws.Range("D24").Value = 0.02929821397056296
$ws.Range("E24").Value = 0.01055645016997686
$ws.Range("D25").Value = 0.0155516254749216
$ws.Range("E25").Value = 0.03921348314606732
$ws.Range("D26").Value = 0.01503090812755905
$ws.Range("E26").Value = 0.02650546384561703
$ws.Range("D27").Value = 0.03070892696308207
$ws.Range("E27").Value = -0.002851380485430099
$ws.Range("D28").Value = 0.03044623845786898
$ws.Range("E28").Value = -0.008774614835220906
$ws.Range("D29").Value = 0.02932325965953007
$ws.Range("E29").Value = -0.0009137136500873
$ws.Range("D30").Value = 0.02868643904083166
$ws.Range("E30").Value = 0.01913341274568192
$ws.Range("D31").Value = 0.03645662135209055
$ws.Range("E31").Value = 0.003339138214759396
$ws.Range("D32").Value = 0.03167211814886342
$ws.Range("E32").Value = 0
$ws.Range("D33").Value = 0.02884409097448504
$ws.Range("E33").Value = 0.01505075253762689
$ws.Range("D34").Value = 0.03207789714065601
$ws.Range("E34").Value = 0.0003873623048056274
$ws.Range("D35").Value = 0.03015966917945469
$ws.Range("E35").Value = 0.0002317497103128119
$ws.Range("D36").Value = 0.03142360278546891
$ws.Range("E36").Value = -0.0001235712079084017
$ws.Range("D37").Value = 0.03499406960798837
$ws.Range("E37").Value = 0.008033732800710069
$ws.Range("D38").Value = 1
$ws.Range("E38").Value = 0.004048276439170717

# Restore sheet protection (content-protected, matching the original state).
$ws.Protect("")

Write-Output "edit complete"
